# Functionality to check the KeyWordDriven FrameWork from the facebook website
# Update the "Enter Password" step on the Login sheet to locate the password
# field with a CSS selector instead of an id locator, and leave the Login
# sheet active/selected (matching the saved view state of the workbook).

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("Login")

# Row 5 = "Enter Password" step: Locator -> cssSelector, LocatorValue -> #pass
$loginSheet.Range("C5").Value = "#pass"
$loginSheet.Range("B5").Value = "cssSelector"

# Make the Login sheet the active sheet/tab, with B5 selected.
$loginSheet.Activate()
[void]$loginSheet.Range("B5").Select()
